$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in column H (row 1), matching the style of the
# other header cells (B1:G1) by copying G1's formatting onto H1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new "Save" column values for the data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
